# Apply the division-problem text replacements described by the diff.
# Each "old" text is unique within the document, so a simple
# Find/Replace (wdReplaceAll) per pair is unambiguous and safe.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "68÷2="; New = "84÷5=" },
    @{ Old = "58÷4="; New = "95÷5=" },
    @{ Old = "26÷4="; New = "53÷6=" },
    @{ Old = "23÷4="; New = "78÷9=" },
    @{ Old = "50÷5="; New = "13÷4=" },
    @{ Old = "24÷8="; New = "64÷9=" },
    @{ Old = "74÷8="; New = "57÷9=" },
    @{ Old = "80÷4="; New = "97÷3=" },
    @{ Old = "91÷6="; New = "95÷5=" },
    @{ Old = "36÷5="; New = "19÷8=" },
    @{ Old = "76÷5="; New = "36÷4=" },
    @{ Old = "38÷8="; New = "34÷7=" },
    @{ Old = "22÷3="; New = "61÷9=" },
    @{ Old = "75÷5="; New = "20÷6=" },
    @{ Old = "19÷6="; New = "77÷3=" },
    @{ Old = "64÷3="; New = "72÷6=" },
    @{ Old = "33÷7="; New = "41÷5=" },
    @{ Old = "36÷9="; New = "68÷8=" },
    @{ Old = "11÷4="; New = "85÷9=" },
    @{ Old = "12÷9="; New = "48÷4=" },
    @{ Old = "47÷8="; New = "39÷4=" },
    @{ Old = "77÷2="; New = "38÷6=" },
    @{ Old = "22÷4="; New = "16÷8=" },
    @{ Old = "62÷8="; New = "25÷9=" },
    @{ Old = "92÷9="; New = "25÷2=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false,
                         $true, 1, $false, $r.New, 2)
}

$d.Save()
